# Adds two new "Mac-Address" test users (Jane Smith in row 31, John Doe in
# row 32) to the bottom of the master-user_detail table, mirroring the
# existing rows' layout/format, and updates the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- name / email columns first -------------------------------------------
# Shared strings must be interned in this precise order (John Doe,
# john.doe@xyz.com, Jane Smith, jane.smith@xyz.com), so populate row 32's
# name/email ahead of row 31's, even though row 31 is written first overall.
#
# Column D (email) carries the same cell style as the existing rows
# (style index 2). A plain `.Value =` on a brand-new cell does not carry
# that formatting, so we first "Insert" the copied, already-styled source
# cell (which shifts nothing, since there is no data below) and only then
# overwrite it with the real value - this keeps the style while landing the
# correct final text/shared-string.
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D32").Insert(-4121) | Out-Null   # xlShiftDown
$excel.CutCopyMode = $false
$ws.Cells.Item(32, 3).Value = "John Doe"
$ws.Cells.Item(32, 4).Value = "john.doe@xyz.com"

$ws.Range("D30").Copy() | Out-Null
$ws.Range("D31").Insert(-4121) | Out-Null   # xlShiftDown
$excel.CutCopyMode = $false
$ws.Cells.Item(31, 3).Value = "Jane Smith"
$ws.Cells.Item(31, 4).Value = "jane.smith@xyz.com"

# --- row 31: Jane Smith -----------------------------------------------------
$ws.Cells.Item(31, 1).Value = 110030
$ws.Cells.Item(31, 2).Value = 9317596768
$ws.Cells.Item(31, 5).Value = 818876432
$ws.Cells.Item(31, 6).Value = "ACT"
$ws.Cells.Item(31, 7).Value = "eng"
$ws.Cells.Item(31, 8).Value = "PWD"
$ws.Cells.Item(31, 9).Value = $true
$ws.Cells.Item(31, 10).Value = "superadmin"
$ws.Cells.Item(31, 11).Value = "now()"
$ws.Range("I31").HorizontalAlignment = -4131   # xlLeft -> reuses style index 1

# --- row 32: John Doe --------------------------------------------------------
$ws.Cells.Item(32, 1).Value = 110031
$ws.Cells.Item(32, 2).Value = 9317596767
$ws.Cells.Item(32, 5).Value = 818876431
$ws.Cells.Item(32, 6).Value = "ACT"
$ws.Cells.Item(32, 7).Value = "eng"
$ws.Cells.Item(32, 8).Value = "PWD"
$ws.Cells.Item(32, 9).Value = $true
$ws.Cells.Item(32, 10).Value = "superadmin"
$ws.Cells.Item(32, 11).Value = "now()"
$ws.Range("I32").HorizontalAlignment = -4131   # xlLeft -> reuses style index 1

# --- view: scroll down a bit and select the new cell near the bottom -------
$ws.Activate() | Out-Null
$ws.Range("E28").Select() | Out-Null
